$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1212.5
$ws.Range("J17").Value = 1212.5
$ws.Range("L17").Value = 3637.5
$ws.Range("N17").Value = -3973.5

$ws.Range("H82").Value = 780.6667
$ws.Range("I82").Value = 780.6667
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2342.0001
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1936.0001
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 780.6667
$ws.Range("I85").Value = 780.6667
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2342.0001
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -938.0001000000002
$ws.Range("N85").ClearContents()

$ws.Range("H98").Value = 3835.8696
$ws.Range("I98").Value = 4370.4736
$ws.Range("J98").Value = 1296.5
$ws.Range("K98").Value = 4370.4736
$ws.Range("L98").Value = 1296.5
$ws.Range("M98").Value = -2872.4736
$ws.Range("N98").Value = -4292.5

$ws.Range("H116").Value = 2212.6316
$ws.Range("I116").Value = 1763.4
$ws.Range("J116").Value = 2711.7778
$ws.Range("K116").Value = 1763.4
$ws.Range("L116").Value = 2711.7778
$ws.Range("M116").Value = 1678.6
$ws.Range("N116").Value = -9595.7778

$ws.Range("H122").Value = 3835.8696
$ws.Range("I122").Value = 4370.4736
$ws.Range("J122").Value = 1296.5
$ws.Range("K122").Value = 13111.4208
$ws.Range("L122").Value = 3889.5
$ws.Range("M122").Value = -10661.4208
$ws.Range("N122").Value = -8789.5

$ws.Range("H129").Value = 707.0303
$ws.Range("I129").Value = 355.8
$ws.Range("J129").Value = 859.73914
$ws.Range("K129").Value = 1067.4
$ws.Range("L129").Value = 2579.21742
$ws.Range("M129").Value = 3932.6
$ws.Range("N129").Value = -12579.21742

$ws.Range("H135").Value = 446.63635
$ws.Range("I135").Value = 216.53847
$ws.Range("J135").Value = 2241.4
$ws.Range("K135").Value = 1948.84623
$ws.Range("L135").Value = 20172.6
$ws.Range("M135").Value = 586.1537700000001
$ws.Range("N135").Value = -25242.6

$ws.Range("H137").Value = 1282.6863
$ws.Range("I137").Value = 888.6177
$ws.Range("J137").Value = 2070.8235
$ws.Range("K137").Value = 2665.8531
$ws.Range("L137").Value = 6212.470499999999
$ws.Range("M137").Value = -115.8531000000003
$ws.Range("N137").Value = -11312.4705

$ws.Range("H138").Value = 1172.5
$ws.Range("I138").Value = 684.6818
$ws.Range("K138").Value = 2054.0454
$ws.Range("M138").Value = 3085.9546

$ws.Range("H141").Value = 651
$ws.Range("I141").Value = 651
$ws.Range("K141").Value = 1953
$ws.Range("M141").Value = 3227

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4101.644
$ws.Range("I32").Value = 3630.5588
$ws.Range("J32").Value = 10508.4
$ws.Range("K32").Value = 3630.5588
$ws.Range("L32").Value = 10508.4
$ws.Range("M32").Value = -3343.5588
$ws.Range("N32").Value = -11082.4

$ws.Range("H61").Value = 47620692
$ws.Range("I61").Value = 62501464
$ws.Range("J61").Value = 2222.8
$ws.Range("K61").Value = 62501464
$ws.Range("L61").Value = 2222.8
$ws.Range("M61").Value = -62501252
$ws.Range("N61").Value = -2646.8

$ws.Range("H74").Value = 1481.4783
$ws.Range("I74").Value = 1055.7894
$ws.Range("K74").Value = 1055.7894
$ws.Range("M74").Value = -181.7893999999999

$ws.Range("H77").Value = 1481.4783
$ws.Range("I77").Value = 1055.7894
$ws.Range("K77").Value = 5278.946999999999
$ws.Range("M77").Value = -910.9469999999992

$ws.Range("H132").Value = 1899.4572
$ws.Range("I132").Value = 1875.2069
$ws.Range("J132").Value = 2016.6666
$ws.Range("K132").Value = 5625.620699999999
$ws.Range("L132").Value = 6049.9998
$ws.Range("M132").Value = -3095.620699999999
$ws.Range("N132").Value = -11109.9998

$ws.Range("H136").Value = 47620692
$ws.Range("I136").Value = 62501464
$ws.Range("J136").Value = 2222.8
$ws.Range("K136").Value = 187504392
$ws.Range("L136").Value = 6668.400000000001
$ws.Range("M136").Value = -187501842
$ws.Range("N136").Value = -11768.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 222
$ws.Range("I11").Value = 222
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 222
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -82
$ws.Range("N11").ClearContents()

$ws.Range("H86").Value = 2881.2778
$ws.Range("I86").Value = 3059.75
$ws.Range("K86").Value = 3059.75
$ws.Range("M86").Value = -1936.75

$ws.Range("H89").Value = 2881.2778
$ws.Range("I89").Value = 3059.75
$ws.Range("K89").Value = 15298.75
$ws.Range("M89").Value = -9682.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1852.1515
$ws.Range("I31").Value = 1838.1111
$ws.Range("J31").Value = 1915.3334
$ws.Range("K31").Value = 1838.1111
$ws.Range("L31").Value = 1915.3334
$ws.Range("M31").Value = -1543.1111
$ws.Range("N31").Value = -2505.3334

$ws.Range("H34").Value = 1852.1515
$ws.Range("I34").Value = 1838.1111
$ws.Range("J34").Value = 1915.3334
$ws.Range("K34").Value = 1838.1111
$ws.Range("L34").Value = 1915.3334
$ws.Range("M34").Value = -1636.1111
$ws.Range("N34").Value = -2319.3334

$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 3000
$ws.Range("K39").Value = 3000
$ws.Range("M39").Value = -2609

$ws.Range("H49").Value = 3000
$ws.Range("I49").Value = 3000
$ws.Range("K49").Value = 3000
$ws.Range("M49").Value = -2818

$ws.Range("H68").Value = 16500
$ws.Range("J68").Value = 16500
$ws.Range("L68").Value = 16500
$ws.Range("N68").Value = -17998

$ws.Range("H71").Value = 16500
$ws.Range("J71").Value = 16500
$ws.Range("L71").Value = 49500
$ws.Range("N71").Value = -56988

$ws.Range("H99").Value = 1688.4667
$ws.Range("I99").Value = 1883.5555
$ws.Range("J99").Value = 1395.8334
$ws.Range("K99").Value = 1883.5555
$ws.Range("L99").Value = 1395.8334
$ws.Range("M99").Value = -385.5554999999999
$ws.Range("N99").Value = -4391.8334

$ws.Range("H126").Value = 1688.4667
$ws.Range("I126").Value = 1883.5555
$ws.Range("J126").Value = 1395.8334
$ws.Range("K126").Value = 5650.666499999999
$ws.Range("L126").Value = 4187.5002
$ws.Range("M126").Value = -3180.666499999999
$ws.Range("N126").Value = -9127.5002

$ws.Range("H134").Value = 12196468
$ws.Range("I134").Value = 1373.1316
$ws.Range("J134").Value = 166667660
$ws.Range("K134").Value = 4119.3948
$ws.Range("L134").Value = 500002980
$ws.Range("M134").Value = -1584.3948
$ws.Range("N134").Value = -500008050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2308.3635
$ws.Range("I5").Value = 2665.7778
$ws.Range("K5").Value = 7997.3334
$ws.Range("M5").Value = -7885.3334

$ws.Range("H39").Value = 2637.32
$ws.Range("J39").Value = 2549.1904
$ws.Range("L39").Value = 7647.5712
$ws.Range("N39").Value = -8235.5712

$ws.Range("H55").Value = 2627.8572
$ws.Range("J55").Value = 2982.5
$ws.Range("L55").Value = 8947.5
$ws.Range("N55").Value = -9301.5

$ws.Range("H113").Value = 729.8182
$ws.Range("J113").Value = 729.8182
$ws.Range("L113").Value = 2189.4546
$ws.Range("N113").Value = -6529.4546

$ws.Range("H131").Value = 17858372
$ws.Range("I131").Value = 83333896
$ws.Range("J131").Value = 1411.7727
$ws.Range("K131").Value = 250001688
$ws.Range("L131").Value = 4235.3181
$ws.Range("M131").Value = -249996648
$ws.Range("N131").Value = -14315.3181

$ws.Range("H135").Value = 2308.3635
$ws.Range("I135").Value = 2665.7778
$ws.Range("K135").Value = 23992.0002
$ws.Range("M135").Value = -21457.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1446.3572
$ws.Range("I113").Value = 1333.6
$ws.Range("K113").Value = 1333.6
$ws.Range("M113").Value = 836.4000000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 899.2857
$ws.Range("I22").Value = 599.0909
$ws.Range("K22").Value = 599.0909
$ws.Range("M22").Value = -304.0909

$ws.Range("H27").Value = 899.2857
$ws.Range("I27").Value = 599.0909
$ws.Range("K27").Value = 599.0909
$ws.Range("M27").Value = -492.0909

$ws.Range("H122").Value = 11371619
$ws.Range("I122").Value = 17865872
$ws.Range("J122").Value = 6675.625
$ws.Range("K122").Value = 53597616
$ws.Range("L122").Value = 20026.875
$ws.Range("M122").Value = -53595166
$ws.Range("N122").Value = -24926.875

$ws.Range("H132").Value = 17937.574
$ws.Range("I132").Value = 1394.725
$ws.Range("J132").Value = 49447.76
$ws.Range("K132").Value = 4184.174999999999
$ws.Range("L132").Value = 148343.28
$ws.Range("M132").Value = -1654.174999999999
$ws.Range("N132").Value = -153403.28

$ws.Range("H136").Value = 4943.6665
$ws.Range("I136").Value = 5525.4165
$ws.Range("J136").Value = 2616.6667
$ws.Range("K136").Value = 16576.2495
$ws.Range("L136").Value = 7850.000100000001
$ws.Range("M136").Value = -14026.2495
$ws.Range("N136").Value = -12950.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H107").Value = 525.4545000000001
$ws.Range("I107").Value = 448
$ws.Range("K107").Value = 1344
$ws.Range("M107").Value = 576

$ws.Range("H122").Value = 8624013
$ws.Range("I122").Value = 10420415
$ws.Range("J122").Value = 1281.8
$ws.Range("K122").Value = 31261245
$ws.Range("L122").Value = 3845.4
$ws.Range("M122").Value = -31258795
$ws.Range("N122").Value = -8745.4

$ws.Range("H132").Value = 3164.963
$ws.Range("I132").Value = 3664.476
$ws.Range("J132").Value = 1416.6666
$ws.Range("K132").Value = 10993.428
$ws.Range("L132").Value = 4249.9998
$ws.Range("M132").Value = -8463.428
$ws.Range("N132").Value = -9309.9998

$ws.Range("H136").Value = 474.9355
$ws.Range("I136").Value = 390.1154
$ws.Range("J136").Value = 916
$ws.Range("K136").Value = 1170.3462
$ws.Range("L136").Value = 2748
$ws.Range("M136").Value = 1379.6538
$ws.Range("N136").Value = -7848
